$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '275.23'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.37%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.56'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-2.67%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.885'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.88%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06346'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.59%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.909'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.21%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.313'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.31%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.267'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '34.83%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8681'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.02%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1537'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '5.11%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.05007'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-2.97%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07411'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.38%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.02970'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-5.38%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09044'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.29%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001572'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.47%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0006326'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.51%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.005966'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.53%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.447'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.04%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.56%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1325'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.05%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.908'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.50%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04351'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.68%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001179'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.41%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004213'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-1.75%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.04%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001678'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.67%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04102'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.84%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007004'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '7.01%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1169'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.73%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-1.39%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01079'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-18.32%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005276'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.20%'
$ws.Range("B46").Value = 'BOLO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.486'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-37.49%'
$ws.Range("B47").Value = 'CoinbaseStockToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.01999'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-11.20%'
